$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3425026666666667
$ws.Range("H2").Value = 1.027508
$ws.Range("I2").Value = 0.00298389995661603
$ws.Range("J2").Value = 0.00298389995661603
$ws.Range("M2").Value = 24.75542533333333
$ws.Range("N2").Value = 74.26627599999999
$ws.Range("O2").Value = 0.7762421087066456
$ws.Range("P2").Value = 0.7762421087066456
$ws.Range("Q2").Value = 8.478799191134222
$ws.Range("R2").Value = 76.309192720208
$ws.Range("S2").Value = 0.002316228794493295
$ws.Range("T2").Value = 0.002316228794493295
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3425026666666667
$ws.Range("H3").Value = 1.027508
$ws.Range("I3").Value = 0.00298389995661603
$ws.Range("J3").Value = 0.00298389995661603
$ws.Range("M3").Value = 3.818542
$ws.Range("N3").Value = 11.455626
$ws.Range("O3").Value = 0.1197358984688377
$ws.Range("P3").Value = 0.1197358984688377
$ws.Range("Q3").Value = 1.307860817778667
$ws.Range("R3").Value = 11.770747360008
$ws.Range("S3").Value = 0.0003572799422465461
$ws.Range("T3").Value = 0.0003572799422465461
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3425026666666667
$ws.Range("H4").Value = 1.027508
$ws.Range("I4").Value = 0.00298389995661603
$ws.Range("J4").Value = 0.00298389995661603
$ws.Range("M4").Value = 3.317404
$ws.Range("N4").Value = 9.952212
$ws.Range("O4").Value = 0.1040219928245168
$ws.Range("P4").Value = 0.1040219928245168
$ws.Range("Q4").Value = 1.136219716410667
$ws.Range("R4").Value = 10.225977447696
$ws.Range("S4").Value = 0.0003103912198761886
$ws.Range("T4").Value = 0.0003103912198761886
$ws.Range("I5").Value = 0.8482396435584867
$ws.Range("J5").Value = 0.8482396435584867
$ws.Range("M5").Value = 24.75542533333333
$ws.Range("N5").Value = 74.26627599999999
$ws.Range("O5").Value = 0.7762421087066456
$ws.Range("P5").Value = 0.7762421087066456
$ws.Range("Q5").Value = 2410.286439981056
$ws.Range("R5").Value = 21692.5779598295
$ws.Range("S5").Value = 0.6584393296044131
$ws.Range("T5").Value = 0.6584393296044131
$ws.Range("I6").Value = 0.8482396435584867
$ws.Range("J6").Value = 0.8482396435584867
$ws.Range("M6").Value = 3.818542
$ws.Range("N6").Value = 11.455626
$ws.Range("O6").Value = 0.1197358984688377
$ws.Range("P6").Value = 0.1197358984688377
$ws.Range("Q6").Value = 371.788401094656
$ws.Range("R6").Value = 3346.095609851904
$ws.Range("S6").Value = 0.101564735838362
$ws.Range("T6").Value = 0.101564735838362
$ws.Range("I7").Value = 0.8482396435584867
$ws.Range("J7").Value = 0.8482396435584867
$ws.Range("M7").Value = 3.317404
$ws.Range("N7").Value = 9.952212
$ws.Range("O7").Value = 0.1040219928245168
$ws.Range("P7").Value = 0.1040219928245168
$ws.Range("Q7").Value = 322.995616899072
$ws.Range("R7").Value = 2906.960552091648
$ws.Range("S7").Value = 0.08823557811571156
$ws.Range("T7").Value = 0.08823557811571157
$ws.Range("G8").Value = 17.07709166666666
$ws.Range("H8").Value = 51.231275
$ws.Range("I8").Value = 0.1487764564848973
$ws.Range("J8").Value = 0.1487764564848973
$ws.Range("M8").Value = 24.75542533333333
$ws.Range("N8").Value = 74.26627599999999
$ws.Range("O8").Value = 0.7762421087066456
$ws.Range("P8").Value = 0.7762421087066456
$ws.Range("Q8").Value = 422.7506676646555
$ws.Range("R8").Value = 3804.756008981899
$ws.Range("S8").Value = 0.1154865503077392
$ws.Range("T8").Value = 0.1154865503077392
$ws.Range("G9").Value = 17.07709166666666
$ws.Range("H9").Value = 51.231275
$ws.Range("I9").Value = 0.1487764564848973
$ws.Range("J9").Value = 0.1487764564848973
$ws.Range("M9").Value = 3.818542
$ws.Range("N9").Value = 11.455626
$ws.Range("O9").Value = 0.1197358984688377
$ws.Range("P9").Value = 0.1197358984688377
$ws.Range("Q9").Value = 65.20959176701666
$ws.Range("R9").Value = 586.88632590315
$ws.Range("S9").Value = 0.01781388268822911
$ws.Range("T9").Value = 0.01781388268822911
$ws.Range("G10").Value = 17.07709166666666
$ws.Range("H10").Value = 51.231275
$ws.Range("I10").Value = 0.1487764564848973
$ws.Range("J10").Value = 0.1487764564848973
$ws.Range("M10").Value = 3.317404
$ws.Range("N10").Value = 9.952212
$ws.Range("O10").Value = 0.1040219928245168
$ws.Range("P10").Value = 0.1040219928245168
$ws.Range("Q10").Value = 56.65161220336665
$ws.Range("R10").Value = 509.8645098302999
$ws.Range("S10").Value = 0.01547602348892902
$ws.Range("T10").Value = 0.01547602348892902